$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 23037.75
$ws.Range("I7").Value = 3002.5
$ws.Range("J7").Value = 29716.166
$ws.Range("K7").Value = 3002.5
$ws.Range("L7").Value = 29716.166
$ws.Range("M7").Value = -2890.5
$ws.Range("N7").Value = -29940.166

$ws.Range("H14").Value = 23037.75
$ws.Range("I14").Value = 3002.5
$ws.Range("J14").Value = 29716.166
$ws.Range("K14").Value = 3002.5
$ws.Range("L14").Value = 29716.166
$ws.Range("M14").Value = -2811.5
$ws.Range("N14").Value = -30098.166

$ws.Range("H62").Value = 1062.5714
$ws.Range("I62").Value = 1170
$ws.Range("J62").Value = 794
$ws.Range("K62").Value = 1170
$ws.Range("L62").Value = 794
$ws.Range("M62").Value = -546
$ws.Range("N62").Value = -2042

$ws.Range("H65").Value = 1062.5714
$ws.Range("I65").Value = 1170
$ws.Range("J65").Value = 794
$ws.Range("K65").Value = 5850
$ws.Range("L65").Value = 3970
$ws.Range("M65").Value = -2730
$ws.Range("N65").Value = -10210

$ws.Range("H137").Value = 1505.7142
$ws.Range("I137").Value = 1440
$ws.Range("J137").Value = 1670
$ws.Range("K137").Value = 4320
$ws.Range("L137").Value = 5010
$ws.Range("M137").Value = -1770
$ws.Range("N137").Value = -10110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4676.5
$ws.Range("I32").Value = 3341.775
$ws.Range("K32").Value = 3341.775
$ws.Range("M32").Value = -3054.775

$ws.Range("H61").Value = 9871.538
$ws.Range("I61").Value = 11312.363
$ws.Range("J61").Value = 1947
$ws.Range("K61").Value = 11312.363
$ws.Range("L61").Value = 1947
$ws.Range("M61").Value = -11100.363
$ws.Range("N61").Value = -2371

$ws.Range("H74").Value = 5418.357
$ws.Range("I74").Value = 5867.25
$ws.Range("K74").Value = 5867.25
$ws.Range("M74").Value = -4993.25

$ws.Range("H77").Value = 5418.357
$ws.Range("I77").Value = 5867.25
$ws.Range("K77").Value = 29336.25
$ws.Range("M77").Value = -24968.25

$ws.Range("H110").Value = 1599.5
$ws.Range("I110").Value = 1199
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 1199
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 846
$ws.Range("N110").Value = -6090

$ws.Range("H136").Value = 9871.538
$ws.Range("I136").Value = 11312.363
$ws.Range("J136").Value = 1947
$ws.Range("K136").Value = 33937.089
$ws.Range("L136").Value = 5841
$ws.Range("M136").Value = -31387.089
$ws.Range("N136").Value = -10941

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2125.4546
$ws.Range("I107").Value = 2233.3333
$ws.Range("J107").Value = 1996
$ws.Range("K107").Value = 2233.3333
$ws.Range("L107").Value = 1996
$ws.Range("M107").Value = -313.3332999999998
$ws.Range("N107").Value = -5836

$ws.Range("H134").Value = 7125.5
$ws.Range("I134").Value = 14306.444
$ws.Range("J134").Value = 2816.9333
$ws.Range("K134").Value = 42919.33199999999
$ws.Range("L134").Value = 8450.7999
$ws.Range("M134").Value = -40384.33199999999
$ws.Range("N134").Value = -13520.7999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 251.86667
$ws.Range("I7").Value = 256.2857
$ws.Range("K7").Value = 256.2857
$ws.Range("M7").Value = -143.2857

$ws.Range("H22").Value = 433.41177
$ws.Range("J22").Value = 507.125
$ws.Range("L22").Value = 507.125
$ws.Range("N22").Value = -1207.125

$ws.Range("H31").Value = 8208.941000000001
$ws.Range("I31").Value = 1733.6923
$ws.Range("K31").Value = 1733.6923
$ws.Range("M31").Value = -1438.6923

$ws.Range("H34").Value = 8208.941000000001
$ws.Range("I34").Value = 1733.6923
$ws.Range("K34").Value = 1733.6923
$ws.Range("M34").Value = -1531.6923

$ws.Range("H58").Value = 1022.4483
$ws.Range("I58").Value = 835.7826
$ws.Range("J58").Value = 1738
$ws.Range("K58").Value = 835.7826
$ws.Range("L58").Value = 1738
$ws.Range("M58").Value = -632.7826
$ws.Range("N58").Value = -2144

$ws.Range("H107").Value = 986.3333
$ws.Range("I107").Value = 986.3333
$ws.Range("K107").Value = 986.3333
$ws.Range("M107").Value = 933.6667

$ws.Range("H132").Value = 3859.8
$ws.Range("I132").Value = 3485.7144
$ws.Range("K132").Value = 10457.1432
$ws.Range("M132").Value = -7927.143199999999

$ws.Range("H134").Value = 3693.5789
$ws.Range("I134").Value = 3979.875
$ws.Range("J134").Value = 2166.6667
$ws.Range("K134").Value = 11939.625
$ws.Range("L134").Value = 6500.000100000001
$ws.Range("M134").Value = -9404.625
$ws.Range("N134").Value = -11570.0001

$ws.Range("H136").Value = 1022.4483
$ws.Range("I136").Value = 835.7826
$ws.Range("J136").Value = 1738
$ws.Range("K136").Value = 2507.3478
$ws.Range("L136").Value = 5214
$ws.Range("M136").Value = 42.65219999999999
$ws.Range("N136").Value = -10314

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6610.143
$ws.Range("I126").Value = 8836.571
$ws.Range("J126").Value = 2157.2856
$ws.Range("K126").Value = 26509.713
$ws.Range("L126").Value = 6471.8568
$ws.Range("M126").Value = -24039.713
$ws.Range("N126").Value = -11411.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 841
$ws.Range("I16").Value = 561.1875
$ws.Range("J16").Value = 2333.3333
$ws.Range("K16").Value = 561.1875
$ws.Range("L16").Value = 2333.3333
$ws.Range("M16").Value = -391.1875
$ws.Range("N16").Value = -2673.3333

$ws.Range("H22").Value = 3833961.5
$ws.Range("I22").Value = 55556050
$ws.Range("J22").Value = 2695.926
$ws.Range("K22").Value = 55556050
$ws.Range("L22").Value = 2695.926
$ws.Range("M22").Value = -55555755
$ws.Range("N22").Value = -3285.926

$ws.Range("H27").Value = 3833961.5
$ws.Range("I27").Value = 55556050
$ws.Range("J27").Value = 2695.926
$ws.Range("K27").Value = 55556050
$ws.Range("L27").Value = 2695.926
$ws.Range("M27").Value = -55555943
$ws.Range("N27").Value = -2909.926

$ws.Range("H132").Value = 14948557
$ws.Range("I132").Value = 19703740
$ws.Range("J132").Value = 3696.8572
$ws.Range("K132").Value = 59111220
$ws.Range("L132").Value = 11090.5716
$ws.Range("M132").Value = -59108690
$ws.Range("N132").Value = -16150.5716

$ws.Range("H136").Value = 24455.334
$ws.Range("I136").Value = 28549.75
$ws.Range("J136").Value = 21179.8
$ws.Range("K136").Value = 85649.25
$ws.Range("L136").Value = 63539.39999999999
$ws.Range("M136").Value = -83099.25
$ws.Range("N136").Value = -68639.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 76923560
$ws.Range("I107").Value = 83333790
$ws.Range("K107").Value = 250001370
$ws.Range("M107").Value = -249999450

$ws.Range("H132").Value = 2031.9259
$ws.Range("I132").Value = 1370
$ws.Range("J132").Value = 2487
$ws.Range("K132").Value = 4110
$ws.Range("L132").Value = 7461
$ws.Range("M132").Value = -1580
$ws.Range("N132").Value = -12521

$ws.Range("H136").Value = 1465.4615
$ws.Range("I136").Value = 850.6667
$ws.Range("K136").Value = 2552.0001
$ws.Range("M136").Value = -2.000100000000202
